$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.179.57"
$ws.Range("E2").Value = "  +3.54%  "

$ws.Range("D3").Value = "1.903.22"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.98"
$ws.Range("E5").Value = "  +3.33%  "

$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4006"
$ws.Range("E8").Value = "  +1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08462"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.66"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.27"
$ws.Range("E12").Value = "  +13.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.430"
$ws.Range("E13").Value = "  +2.76%  "

$ws.Range("D14").Value = "1.911.75"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.350"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001112"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.35"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.995"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").Value = "30.168.85"
$ws.Range("E23").Value = "  +3.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.206"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").Value = "2.129.91"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.57"
$ws.Range("E27").Value = "  +3.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.09"
$ws.Range("E28").Value = "  +2.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.385"
$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.71"
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("E31").Value = "  +3.56%  "

$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.050"
$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.676"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02492"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06561"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2200"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.188"
$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.93"
$ws.Range("E40").Value = "  +5.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.784"
$ws.Range("E41").Value = "  -2.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6506"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.231"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6122"
$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.716"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.243"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.08"
$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.162"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.18"
$ws.Range("E51").Value = "  +2.09%  "

